# Fix to intimate partner calcs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Table 1 (rows 5-14): turn individual D formulas into one shared-formula fill ---
$ws.Range("D5:D14").Formula = "=C5/SUM(B5:C5)"

# --- Column F width (bestFit-ish) ---
$ws.Columns.Item(6).ColumnWidth = 19.7142857142857

# --- New "Of people who had sex" section header row (row 20) ---
# Shared-string insertion order in the target file is:
#   "Of people who had sex", "SD Prop", "SD", "Difference calc"
# so touch the "SD Prop" cells before "SD"/"Difference calc" to match.
$ws.Range("F20").Value = "Of people who had sex"
$ws.Range("C30").Value = "SD Prop"
$ws.Range("G20").Value = "SD"
$ws.Range("K20").Value = "Difference calc"
$ws.Range("J20").Value = "Year"

# --- Table 2 (rows 21-26): regular partner ---
$ws.Range("D21:D26").Formula = "=C21/SUM(B21:C21)"

$ws.Range("F21").Value = "Only sex with a regular partner"
$ws.Range("F23").Value = "Only sex with non-regular partner"
$ws.Range("F25").Value = "Sex with both regular and non-regular partner"

$ws.Range("F22").Formula = "=F12/F4"
$ws.Range("F22").NumberFormat = "0.0%"
$ws.Range("F24").Formula = "=F14/F6"
$ws.Range("F24").NumberFormat = "0.0%"
$ws.Range("F26").Formula = "=F16/F4"
$ws.Range("F26").NumberFormat = "0.0%"

$ws.Range("G22").Formula = "=D30*F22"
$ws.Range("G24").Formula = "=F24*D42"
$ws.Range("G26").Formula = "=F26*K30"

$ws.Range("J21").Value = 2005
$ws.Range("J22").Formula = "=J21+1"
$ws.Range("J23:J26").Formula = "=J22+1"

# K21 stands alone, K22:K25 are one shared-formula fill, K26 stands alone again
# (mirrors the target file's own formula grouping exactly)
$ws.Range("K21").Formula = '=D21+D33-$D$16'
$ws.Range("K22:K25").Formula = '=D22+D34-$D$16'
$ws.Range("K26").Formula = '=D26+D38-$D$16'
$ws.Range("K21").NumberFormat = "0.0%"
$ws.Range("K22").NumberFormat = "0.0%"
$ws.Range("K23").NumberFormat = "0.0%"
$ws.Range("K24").NumberFormat = "0.0%"
$ws.Range("K25").NumberFormat = "0.0%"
$ws.Range("K26").NumberFormat = "0.0%"

$ws.Range("K27").NumberFormat = "0.0%"

# --- Mean / Standard deviation summary rows ---
$ws.Range("J28").Value = "Mean"
$ws.Range("K28").Formula = "=AVERAGE(K21:K26)"
$ws.Range("K28").NumberFormat = "0.0%"

$ws.Range("J29").Value = "Standard deviation"
$ws.Range("K29").Formula = "=STDEV(K21:K26)"

# --- New SD Prop row (row 30) ---
$ws.Range("D30").Formula = "=D29/D28"

$ws.Range("J30").Value = "SD Prop"
$ws.Range("K30").Formula = "=K29/K28"

$ws.Range("K31").NumberFormat = "0.0%"
$ws.Range("K32").NumberFormat = "0.0%"

# --- Table 3 (rows 33-38): non-regular partner ---
$ws.Range("D33:D38").Formula = "=C33/SUM(B33:C33)"

# --- New SD Prop row (row 42) ---
$ws.Range("C42").Value = "SD Prop"
$ws.Range("D42").Formula = "=D41/D40"

# --- Selection matches final authored state ---
$ws.Range("K30").Select()

$wb.Application.Calculate()
